# ScoreHunter Supports the Cer in the [M+H]+
# Rename the two "[MG(FAx)-H2O+H]+" labels to the new "FAx_[MG-H2O+H]+" form
# used for the other lipid-group rows (FA1/FA2 prefix convention), then
# resize column A so the longer header text is still fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "FA1_[MG-H2O+H]+"
$ws.Range("A5").Value = "FA2_[MG-H2O+H]+"

$ws.Columns.Item(1).AutoFit() | Out-Null
